# Case and Fatality Demographics Data Updated
# - Reorders the "Fatalities by Gender" tab ahead of "Fatalities by Age Group"
# - Refreshes the underlying case counts on the three Fatalities sheets
#   (Age Group, Gender, Race-Ethnicity) with the latest figures; the
#   percentage columns are formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Reorder worksheet tabs: "Fatalities by Gender" moves before
#    "Fatalities by Age Group".
# ---------------------------------------------------------------------
$wsFatGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsFatAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsFatGender.Move($wsFatAge)

# ---------------------------------------------------------------------
# 2) Fatalities by Gender - updated counts
# ---------------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsGender.Range("B2").Value = 29517
$wsGender.Range("B3").Value = 41055
$wsGender.Range("B4").Value = 1
$wsGender.Range("B5").Value = 70573

# ---------------------------------------------------------------------
# 3) Fatalities by Age Group - updated counts
# ---------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsAge.Range("B2").Value = 14
$wsAge.Range("B3").Value = 20
$wsAge.Range("B4").Value = 76
$wsAge.Range("B5").Value = 617
$wsAge.Range("B6").Value = 2002
$wsAge.Range("B7").Value = 4927
$wsAge.Range("B8").Value = 9379
$wsAge.Range("B9").Value = 7151
$wsAge.Range("B10").Value = 8485
$wsAge.Range("B11").Value = 9041
$wsAge.Range("B12").Value = 8570
$wsAge.Range("B13").Value = 20291
$wsAge.Range("B14").Value = 0
$wsAge.Range("B15").Value = 70573

# ---------------------------------------------------------------------
# 4) Fatalities by Race-Ethnicity - updated counts (B8 is a SUM formula
#    and recalculates on its own)
# ---------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsRace.Range("B2").Value = 1307
$wsRace.Range("B3").Value = 7423
$wsRace.Range("B4").Value = 30863
$wsRace.Range("B5").Value = 419
$wsRace.Range("B6").Value = 30517
$wsRace.Range("B7").Value = 44
